# Updates cryptos list prices/volume percentages (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '56.892.35' },
    @{ Cell = 'E2'; Value = '  +0.42%  ' },
    @{ Cell = 'D3'; Value = '2.432.63' },
    @{ Cell = 'E3'; Value = '  -2.34%  ' },
    @{ Cell = 'E4'; Value = '  -0.04%  ' },
    @{ Cell = 'D5'; Value = '489.10' },
    @{ Cell = 'E5'; Value = '  -0.74%  ' },
    @{ Cell = 'D6'; Value = '154.00' },
    @{ Cell = 'E6'; Value = '  +0.88%  ' },
    @{ Cell = 'D7'; Value = '0.999' },
    @{ Cell = 'E7'; Value = '  +0.16%  ' },
    @{ Cell = 'D8'; Value = '0.605' },
    @{ Cell = 'E8'; Value = '  +18.08%  ' },
    @{ Cell = 'D9'; Value = '2.427.98' },
    @{ Cell = 'E9'; Value = '  -2.96%  ' },
    @{ Cell = 'D10'; Value = '6.29' },
    @{ Cell = 'E10'; Value = '  +9.21%  ' },
    @{ Cell = 'D11'; Value = '0.0998' },
    @{ Cell = 'E11'; Value = '  +1.24%  ' },
    @{ Cell = 'D12'; Value = '0.335' },
    @{ Cell = 'E12'; Value = '  +0.22%  ' },
    @{ Cell = 'E13'; Value = '  +1.41%  ' },
    @{ Cell = 'D14'; Value = '2.847.20' },
    @{ Cell = 'E14'; Value = '  -2.56%  ' },
    @{ Cell = 'D15'; Value = '56.982.56' },
    @{ Cell = 'E15'; Value = '  +0.25%  ' },
    @{ Cell = 'E16'; Value = '  -2.90%  ' },
    @{ Cell = 'E17'; Value = '  -2.12%  ' },
    @{ Cell = 'D18'; Value = '2.428.05' },
    @{ Cell = 'E18'; Value = '  -2.78%  ' },
    @{ Cell = 'D19'; Value = '4.72' },
    @{ Cell = 'E19'; Value = '  +3.65%  ' },
    @{ Cell = 'D20'; Value = '323.89' },
    @{ Cell = 'E21'; Value = '  -3.01%  ' },
    @{ Cell = 'D22'; Value = '5.99' },
    @{ Cell = 'E22'; Value = '  +1.64%  ' },
    @{ Cell = 'D23'; Value = '1.00' },
    @{ Cell = 'E23'; Value = '  +0.07%  ' },
    @{ Cell = 'D24'; Value = '58.20' },
    @{ Cell = 'E24'; Value = '  -0.84%  ' },
    @{ Cell = 'E25'; Value = '  -0.59%  ' },
    @{ Cell = 'D26'; Value = '0.998' },
    @{ Cell = 'E26'; Value = '  -0.11%  ' },
    @{ Cell = 'D27'; Value = '0.161' },
    @{ Cell = 'E27'; Value = '  -0.52%  ' },
    @{ Cell = 'D28'; Value = '2.525.58' },
    @{ Cell = 'E28'; Value = '  -2.94%  ' },
    @{ Cell = 'D29'; Value = '7.32' },
    @{ Cell = 'E29'; Value = '  -3.57%  ' },
    @{ Cell = 'D30'; Value = '0.0₃0784' },
    @{ Cell = 'E30'; Value = '  -3.01%  ' },
    @{ Cell = 'D31'; Value = '1.00' },
    @{ Cell = 'E31'; Value = '  +0.20%  ' },
    @{ Cell = 'D32'; Value = '150.45' },
    @{ Cell = 'E32'; Value = '  -0.24%  ' },
    @{ Cell = 'E33'; Value = '  +1.26%  ' },
    @{ Cell = 'D34'; Value = '1.53' },
    @{ Cell = 'E34'; Value = '  +0.23%  ' },
    @{ Cell = 'E36'; Value = '  -0.43%  ' },
    @{ Cell = 'E37'; Value = '  -1.34%  ' },
    @{ Cell = 'D38'; Value = '0.848' },
    @{ Cell = 'E38'; Value = '  -2.65%  ' },
    @{ Cell = 'E39'; Value = '  +9.26%  ' },
    @{ Cell = 'D40'; Value = '34.15' },
    @{ Cell = 'E40'; Value = '  +0.17%  ' },
    @{ Cell = 'E42'; Value = '  -1.60%  ' },
    @{ Cell = 'D43'; Value = '0.997' },
    @{ Cell = 'E43'; Value = '  +0.15%  ' },
    @{ Cell = 'E44'; Value = '  -3.57%  ' },
    @{ Cell = 'D45'; Value = '267.79' },
    @{ Cell = 'E45'; Value = '  -0.67%  ' },
    @{ Cell = 'E46'; Value = '  -5.71%  ' },
    @{ Cell = 'D48'; Value = '0.0228' },
    @{ Cell = 'E48'; Value = '  -0.50%  ' },
    @{ Cell = 'D49'; Value = '4.57' },
    @{ Cell = 'E49'; Value = '  -6.74%  ' },
    @{ Cell = 'D50'; Value = '17.46' },
    @{ Cell = 'E50'; Value = '  -1.75%  ' },
    @{ Cell = 'D51'; Value = '1.875.53' },
    @{ Cell = 'E51'; Value = '  -0.74%  ' }
)

foreach ($u in $updates) {
    $cellRef = $u.Cell
    $val = $u.Value
    $col = ($cellRef -replace '[0-9]+$', '')
    $range = $ws.Range($cellRef)
    if ($col -eq 'D') {
        # Force text storage so numeric-looking strings (e.g. "1.00", "489.10")
        # aren't silently coerced into real numbers by Excel's smart entry.
        $range.NumberFormat = "@"
        $range.Value = $val
        $range.Style = "Normal"
    } else {
        $range.Value = $val
    }
}
